# Apply green highlighting to the "Create post et game" requirement block.
# Word COM maps wdBrightGreen (4) to <w:highlight w:val="green"/>.
$wdBrightGreen = 4

$d = $word.ActiveDocument

# Paragraphs that must be highlighted in their entirety (text + paragraph mark).
# Using Range.Font.HighlightColorIndex (rather than Range.HighlightColorIndex)
# correctly stamps both the run(s) <w:rPr> and the paragraph-mark <w:pPr><w:rPr>.
$fullyHighlighted = @(70, 71, 72, 73, 74, 75, 76, 77, 78, 80, 81, 82, 83, 84, 86, 90, 91, 92, 94, 96, 98, 100, 109, 110)

foreach ($idx in $fullyHighlighted) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Font.HighlightColorIndex = $wdBrightGreen
}

# Paragraph 88 ("Les accès doivent être clairement définis\u00a0:") only has its
# first 41 characters ("Les accès doivent être clairement définis") highlighted;
# the trailing non-breaking space + colon stay unhighlighted, which forces the
# run to split in two.
$p88 = $d.Paragraphs.Item(88)
$r88 = $p88.Range
$head = $r88.Duplicate
$head.SetRange($r88.Start, $r88.Start + 41)
$head.Font.HighlightColorIndex = $wdBrightGreen

Write-Host "Highlighting applied."
